$d = $word.ActiveDocument

# Update the header date
$d.Content.Find.Execute("2025-04-15 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-22 Tuesday", 2)

# Update the multiplication table entries, cell by cell (scoped replace avoids
# collisions between values that are both an "old" text in one cell and a
# "new" text in another cell).
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Old = "16×72=1152"; New = "39×13=507" },
    @{ Row = 1;  Col = 2; Old = "15×62=930";  New = "63×90=5670" },
    @{ Row = 1;  Col = 3; Old = "80×90=7200"; New = "91×25=2275" },
    @{ Row = 1;  Col = 4; Old = "41×22=902";  New = "89×69=6141" },
    @{ Row = 1;  Col = 5; Old = "55×65=3575"; New = "87×69=6003" },

    @{ Row = 5;  Col = 1; Old = "25×88=2200"; New = "72×88=6336" },
    @{ Row = 5;  Col = 2; Old = "14×36=504";  New = "85×70=5950" },
    @{ Row = 5;  Col = 3; Old = "91×96=8736"; New = "98×52=5096" },
    @{ Row = 5;  Col = 4; Old = "50×93=4650"; New = "18×11=198" },
    @{ Row = 5;  Col = 5; Old = "71×68=4828"; New = "61×65=3965" },

    @{ Row = 10; Col = 1; Old = "52×44=2288"; New = "62×59=3658" },
    @{ Row = 10; Col = 2; Old = "36×12=432";  New = "52×88=4576" },
    @{ Row = 10; Col = 3; Old = "81×38=3078"; New = "48×91=4368" },
    @{ Row = 10; Col = 4; Old = "19×69=1311"; New = "48×13=624" },

    # NOTE: row 15 col 3's original value (99×89=8811) is the same text that
    # row 10 col 5 is about to be rewritten to. This cell must be updated
    # first so that the "99×89=8811" search text is still unique in the
    # document when it is looked up (the Find used here resolves to the
    # first match of the search text in the whole document, not strictly
    # within the owning cell's range).
    @{ Row = 15; Col = 3; Old = "99×89=8811"; New = "51×66=3366" },

    @{ Row = 10; Col = 5; Old = "30×39=1170"; New = "99×89=8811" },

    @{ Row = 15; Col = 1; Old = "60×40=2400"; New = "90×22=1980" },
    @{ Row = 15; Col = 2; Old = "13×35=455";  New = "65×29=1885" },
    @{ Row = 15; Col = 4; Old = "93×46=4278"; New = "73×74=5402" },
    @{ Row = 15; Col = 5; Old = "83×93=7719"; New = "52×16=832" },

    @{ Row = 20; Col = 1; Old = "42×20=840";  New = "73×33=2409" },
    @{ Row = 20; Col = 2; Old = "33×30=990";  New = "90×86=7740" },
    @{ Row = 20; Col = 3; Old = "80×22=1760"; New = "66×64=4224" },
    @{ Row = 20; Col = 4; Old = "41×21=861";  New = "25×26=650" },
    @{ Row = 20; Col = 5; Old = "23×59=1357"; New = "39×31=1209" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Find.Execute($u.Old, $true, $false, $false, $false, $false, $true, 1, $false, $u.New, 2)
}
